$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Convert the F9:F14 / G9:G14 / K9:K14 formulas into shared formulas
# (they keep the same values, only the underlying XML representation changes
# to a shared-formula group, matching rows 15+ which already use shared formulas).
$ws.Range("F9:F14").Formula = "=SUM(D9:E9)"
$ws.Range("G9:G14").Formula = "=C9/F9"
$ws.Range("K9:K14").Formula = "=I9/(J9+I9)"

# Add new performance-test rows 37 (no sort) and 38 (sort)
$ws.Range("A37").Value = "no sort"
$ws.Range("B37").Value = 0.5
$ws.Range("C37").Value = 549997283
$ws.Range("D37").Value = 68409785
$ws.Range("E37").Value = 65221693
$ws.Range("F37:F38").Formula = "=SUM(D37:E37)"
$ws.Range("G37:G38").Formula = "=C37/F37"
$ws.Range("K37:K38").Formula = "=I37/(J37+I37)"

$ws.Range("A38").Value = "sort"
$ws.Range("B38").Value = 0.5

# Update the active selection to match the author's final cursor position
$ws.Range("F43").Select() | Out-Null
